$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = 44489
$ws.Range("M2").Value2 = 160
$ws.Range("N2").Value2 = 9500
$ws.Range("O2").Value2 = 10000
$ws.Range("P2").Value2 = 9750
$ws.Range("S2").Value2 = 4875
$ws.Range("D3").Value2 = 44455
$ws.Range("M3").Value2 = 200
$ws.Range("N3").Value2 = 12000
$ws.Range("O3").Value2 = 13000
$ws.Range("P3").Value2 = 12500
$ws.Range("S3").Value2 = 6250
$ws.Range("D4").Value2 = 44874
$ws.Range("M4").Value2 = 300
$ws.Range("N4").Value2 = 7500
$ws.Range("O4").Value2 = 8000
$ws.Range("P4").Value2 = 7750
$ws.Range("S4").Value2 = 3875
$ws.Range("D5").Value2 = 44454
$ws.Range("M5").Value2 = 160
$ws.Range("N5").Value2 = 12000
$ws.Range("O5").Value2 = 13000
$ws.Range("P5").Value2 = 12500
$ws.Range("S5").Value2 = 6250
$ws.Range("D6").Value2 = 44819
$ws.Range("M6").Value2 = 240
$ws.Range("N6").Value2 = 11000
$ws.Range("O6").Value2 = 12000
$ws.Range("P6").Value2 = 11500
$ws.Range("S6").Value2 = 5750
$ws.Range("D7").Value2 = 44461
$ws.Range("M7").Value2 = 200
$ws.Range("N7").Value2 = 11000
$ws.Range("O7").Value2 = 12000
$ws.Range("P7").Value2 = 11500
$ws.Range("S7").Value2 = 5750
$ws.Range("D8").Value2 = 44482
$ws.Range("M8").Value2 = 240
$ws.Range("N8").Value2 = 10000
$ws.Range("O8").Value2 = 11000
$ws.Range("P8").Value2 = 10500
$ws.Range("S8").Value2 = 5250
$ws.Range("D9").Value2 = 44490
$ws.Range("M9").Value2 = 400
$ws.Range("N9").Value2 = 9500
$ws.Range("O9").Value2 = 10000
$ws.Range("P9").Value2 = 9750
$ws.Range("S9").Value2 = 4875
$ws.Range("D10").Value2 = 44475
$ws.Range("M10").Value2 = 240
$ws.Range("N10").Value2 = 11000
$ws.Range("O10").Value2 = 12000
$ws.Range("P10").Value2 = 11500
$ws.Range("S10").Value2 = 5750
$ws.Range("D11").Value2 = 44882
$ws.Range("M11").Value2 = 440
$ws.Range("N11").Value2 = 6000
$ws.Range("O11").Value2 = 7000
$ws.Range("P11").Value2 = 6500
$ws.Range("S11").Value2 = 3250
$ws.Range("D12").Value2 = 44818
$ws.Range("M12").Value2 = 200
$ws.Range("N12").Value2 = 11000
$ws.Range("O12").Value2 = 12000
$ws.Range("P12").Value2 = 11500
$ws.Range("S12").Value2 = 5750
$ws.Range("D14").Value2 = 44895
$ws.Range("M14").Value2 = 240
$ws.Range("N14").Value2 = 3000
$ws.Range("O14").Value2 = 3500
$ws.Range("P14").Value2 = 3250
$ws.Range("S14").Value2 = 1625
$ws.Range("D15").Value2 = 44889
$ws.Range("M15").Value2 = 460
$ws.Range("N15").Value2 = 3500
$ws.Range("O15").Value2 = 4000
$ws.Range("P15").Value2 = 3750
$ws.Range("S15").Value2 = 1875
$ws.Range("D16").Value2 = 44497
$ws.Range("M16").Value2 = 500
$ws.Range("N16").Value2 = 9000
$ws.Range("O16").Value2 = 10000
$ws.Range("P16").Value2 = 9500
$ws.Range("S16").Value2 = 4750
$ws.Range("D17").Value2 = 44875
$ws.Range("M17").Value2 = 400
$ws.Range("N17").Value2 = 7000
$ws.Range("O17").Value2 = 7500
$ws.Range("P17").Value2 = 7250
$ws.Range("S17").Value2 = 3625
$ws.Range("D18").Value2 = 44881
$ws.Range("M18").Value2 = 440
$ws.Range("N18").Value2 = 6000
$ws.Range("O18").Value2 = 7000
$ws.Range("P18").Value2 = 6500
$ws.Range("S18").Value2 = 3250
